# Edit script for LINEA 141 horarios workbook update
# Commit: Horarios actualizados Linea 141 - 593
$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")

# Header updates
$ws1.Range("A2").Value = "Última actualización: 14:53:29"
$ws1.Range("A3").Value = "Total filas: 212"

# Data row updates (re-sorted / new scraped rows)
$ws1.Cells.Item(124, 1).Value = "10:36:50"
$ws1.Cells.Item(124, 2).Value = "12:10"
$ws1.Cells.Item(124, 3).Value = "15_ABASTO"
$ws1.Cells.Item(124, 4).Value = 94
$ws1.Cells.Item(124, 5).Value = "LP1912"

$ws1.Cells.Item(125, 1).Value = "10:36:50"
$ws1.Cells.Item(125, 2).Value = "12:10"
$ws1.Cells.Item(125, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(125, 4).Value = 94
$ws1.Cells.Item(125, 5).Value = "LP1912"

$ws1.Cells.Item(137, 1).Value = "10:36:50"
$ws1.Cells.Item(137, 2).Value = "12:34"
$ws1.Cells.Item(137, 3).Value = "15_ABASTO"
$ws1.Cells.Item(137, 4).Value = 118
$ws1.Cells.Item(137, 5).Value = "LP1912"

$ws1.Cells.Item(138, 1).Value = "11:46:32"
$ws1.Cells.Item(138, 2).Value = "12:34"
$ws1.Cells.Item(138, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(138, 4).Value = 48
$ws1.Cells.Item(138, 5).Value = "LP1912"

$ws1.Cells.Item(190, 1).Value = "13:55:43"
$ws1.Cells.Item(190, 2).Value = "14:51"
$ws1.Cells.Item(190, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(190, 4).Value = 56
$ws1.Cells.Item(190, 5).Value = "LP1912"

$ws1.Cells.Item(191, 1).Value = "13:41:21"
$ws1.Cells.Item(191, 2).Value = "14:51"
$ws1.Cells.Item(191, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(191, 4).Value = 70
$ws1.Cells.Item(191, 5).Value = "LP1912"

$ws1.Cells.Item(204, 1).Value = "14:32:44"
$ws1.Cells.Item(204, 2).Value = "15:53"
$ws1.Cells.Item(204, 3).Value = "10_OLMOS"
$ws1.Cells.Item(204, 4).Value = 81
$ws1.Cells.Item(204, 5).Value = "LP1912"

$ws1.Cells.Item(205, 1).Value = "13:55:43"
$ws1.Cells.Item(205, 2).Value = "15:53"
$ws1.Cells.Item(205, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(205, 4).Value = 118
$ws1.Cells.Item(205, 5).Value = "LP1912"

$ws1.Cells.Item(206, 1).Value = "13:55:43"
$ws1.Cells.Item(206, 2).Value = "15:53"
$ws1.Cells.Item(206, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(206, 4).Value = 118
$ws1.Cells.Item(206, 5).Value = "LP1912"

$ws1.Cells.Item(209, 1).Value = "14:53:29"
$ws1.Cells.Item(209, 2).Value = "16:02"
$ws1.Cells.Item(209, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(209, 4).Value = 69
$ws1.Cells.Item(209, 5).Value = "LP1912"

$ws1.Cells.Item(210, 1).Value = "14:11:28"
$ws1.Cells.Item(210, 2).Value = "16:05"
$ws1.Cells.Item(210, 3).Value = "14_ABASTO"
$ws1.Cells.Item(210, 4).Value = 114
$ws1.Cells.Item(210, 5).Value = "LP1912"

$ws1.Cells.Item(211, 1).Value = "14:46:12"
$ws1.Cells.Item(211, 2).Value = "16:06"
$ws1.Cells.Item(211, 3).Value = "14_ABASTO"
$ws1.Cells.Item(211, 4).Value = 80
$ws1.Cells.Item(211, 5).Value = "LP1912"

$ws1.Cells.Item(212, 1).Value = "14:32:44"
$ws1.Cells.Item(212, 2).Value = "16:14"
$ws1.Cells.Item(212, 3).Value = "17_ROMERO"
$ws1.Cells.Item(212, 4).Value = 102
$ws1.Cells.Item(212, 5).Value = "LP1912"

$ws1.Cells.Item(213, 1).Value = "14:46:12"
$ws1.Cells.Item(213, 2).Value = "16:17"
$ws1.Cells.Item(213, 3).Value = "10_OLMOS"
$ws1.Cells.Item(213, 4).Value = 91
$ws1.Cells.Item(213, 5).Value = "LP1912"

$ws1.Cells.Item(214, 1).Value = "14:32:44"
$ws1.Cells.Item(214, 2).Value = "16:21"
$ws1.Cells.Item(214, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(214, 4).Value = 109
$ws1.Cells.Item(214, 5).Value = "LP1912"

$ws1.Cells.Item(215, 1).Value = "14:46:12"
$ws1.Cells.Item(215, 2).Value = "16:34"
$ws1.Cells.Item(215, 3).Value = "83_ALUAR"
$ws1.Cells.Item(215, 4).Value = 108
$ws1.Cells.Item(215, 5).Value = "LP1912"

$ws1.Cells.Item(216, 1).Value = "14:46:12"
$ws1.Cells.Item(216, 2).Value = "16:41"
$ws1.Cells.Item(216, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(216, 4).Value = 115
$ws1.Cells.Item(216, 5).Value = "LP1912"

$ws1.Cells.Item(217, 1).Value = "14:53:29"
$ws1.Cells.Item(217, 2).Value = "16:46"
$ws1.Cells.Item(217, 3).Value = "17_ROMERO"
$ws1.Cells.Item(217, 4).Value = 113
$ws1.Cells.Item(217, 5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 14:53:29"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 14:53:29"
